# Adds 2022-Q3 data.
#
# Final layout (tab order):
#   1. "总计"     - summary sheet, gains a new 2022-Q3 row, old rows shift down
#   2. "2022-Q3"  - was "2022-Q2", physical sheet keeps its place but gets
#                   brand-new Q3 fund data
#   3. "2022-Q2"  - was "2022-Q1" (physical sheet repurposed); receives the
#                   data that used to live in the "2022-Q2" sheet (copied
#                   before it gets overwritten with Q3 numbers)
#   4. "2022-Q1"  - brand-new physical sheet, a straight duplicate of the
#                   original "2022-Q1" sheet (data untouched)

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")
$q1 = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------------
# 1. Duplicate "2022-Q1" so its original data/format survives unchanged in a
#    brand new physical sheet placed right after it. This new sheet will
#    keep the name "2022-Q1".
# ---------------------------------------------------------------------------
$q1.Copy($null, $q1)
$q1New = $wb.Worksheets.Item(4)
$q1New.Name = "2022-Q1-newtmp"

# ---------------------------------------------------------------------------
# 2. The old "2022-Q1" sheet becomes the new "2022-Q2" sheet: copy over the
#    (still untouched) "2022-Q2" sheet's values + formatting, then fix up
#    its page margins to match what "2022-Q2" used.
# ---------------------------------------------------------------------------
$q2.Range("B1:H8").Copy($q1.Range("B1"))
$q2.Range("A2:A8").Copy($q1.Range("A2"))

$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

$q1.Name = "2022-Q2-newtmp"

# ---------------------------------------------------------------------------
# 3. The old "2022-Q2" sheet becomes the new "2022-Q3" sheet: overwrite the
#    fund data in place with the Q3 figures (names/codes for rows 4-8 stay
#    put, only rows 2-3 swap their fund between each other).
# ---------------------------------------------------------------------------
$q3data = @(
    @("014283", "华夏北交所创新中小企业精选两年定开混合", "3.47", "71.59", "5.95", "0.2065", 4),
    @("014277", "万家北交所慧选两年定期开放混合A",         "3.56", "93.97", "5.54", "0.1972", 7),
    @("014279", "汇添富北交所创新精选两年定开混合A",       "3.20", "93.27", "4.76", "0.1523", 4),
    @("014273", "广发北交所精选两年定开混合A",             "3.37", "64.25", "4.11", "0.1385", 8),
    @("014274", "广发北交所精选两年定开混合C",             "0.85", "64.25", "4.11", "0.0349", 8),
    @("014278", "万家北交所慧选两年定期开放混合C",         "0.49", "93.97", "5.54", "0.0271", 7),
    @("014280", "汇添富北交所创新精选两年定开混合C",       "0.51", "93.27", "4.76", "0.0243", 4)
)

for ($i = 0; $i -lt $q3data.Count; $i++) {
    $r = $i + 2
    $row = $q3data[$i]
    $q2.Cells.Item($r, 2).Value = $row[0]
    $q2.Cells.Item($r, 3).Value = $row[1]
    $q2.Cells.Item($r, 4).Value = $row[2]
    $q2.Cells.Item($r, 5).Value = $row[3]
    $q2.Cells.Item($r, 6).Value = $row[4]
    $q2.Cells.Item($r, 7).Value = $row[5]
    $q2.Cells.Item($r, 8).Value = $row[6]
}

$q2.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 4. Finish renaming the repurposed sheets now that "2022-Q3" freed up its
#    old name.
# ---------------------------------------------------------------------------
$q1.Name = "2022-Q2"
$q1New.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 5. Update the "总计" summary sheet: insert the new 2022-Q3 row, shift the
#    existing 2022-Q2 / 2022-Q1 rows down.
# ---------------------------------------------------------------------------
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q1"
$total.Cells.Item(4, 3).Value = 3
$total.Cells.Item(4, 4).Value = 0.44

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 7
$total.Cells.Item(3, 4).Value = 1.05

$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.78
